# DataTypes.xlsx: add a "Large Double Number:" row right after "Double Number:"
# (row 12), pushing every following row down by one. This mirrors the
# ClosedXML fix that stopped relying on Convert.ToDecimal (which throws for
# doubles outside decimal's range) by exercising a double value close to
# Double.MaxValue (9.999E+307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12, shifting rows 12-41 down to 13-42.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row with the label/value pair.
$ws.Range("B12").Value = "Large Double Number:"
$ws.Range("C12").Value = [double]"9.999E+307"
